$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.495
$ws.Range("A3").Value = -21.557
$ws.Range("B5").Value = 6.234
$ws.Range("C5").Value = -12.404
$ws.Range("C9").Value = -12.057
$ws.Range("C11").Value = -12.495
$ws.Range("A14").Value = -20.891
$ws.Range("A16").Value = -21.363
$ws.Range("B16").Value = 6.008
$ws.Range("C17").Value = -12.048
$ws.Range("A21").Value = -21.898
$ws.Range("C21").Value = -12.984
$ws.Range("A23").Value = -21.709
$ws.Range("A25").Value = -22.078

$wb.Save()
